$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Relay #1 (row 2): nominal current raised from 400 to 1200
$ws.Range("G2").Value = 1200

# Relay #2 (row 3): newly added, mirrors relay #1's settings (A3 is its index)
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 13800
$ws.Range("C3").Value = 5
$ws.Range("D3").Value = 12
$ws.Range("E3").Value = 0.5
$ws.Range("F3").Value = 0.5
$ws.Range("G3").Value = 1200
$ws.Range("H3").Value = 1
$ws.Range("I3").Value = 20
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0.8
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1
$ws.Range("N3").Value = 0.9
$ws.Range("O3").Value = 1.1000000000000001
$ws.Range("P3").Value = 20
$ws.Range("Q3").Value = 1
$ws.Range("R3").Value = 0.7
$ws.Range("S3").Value = 2
$ws.Range("T3").Value = 3
$ws.Range("U3").Value = 0.9
$ws.Range("V3").Value = 1.1000000000000001

# Move the active selection to G4 (matches the post-edit sheetView selection)
$ws.Range("G4").Select()
